# Update of "Bosnia Herzegovina Premier Liga" base, 01-04-2024
# - Fixes 6 rows where Home/Away data had been swapped between two
#   consecutive fixtures (columns B:AC swapped, A/row index untouched).
# - Appends one new fixture row (row 145).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# Pairs of rows whose match data (everything except the running index in
# column A) had been mixed up and need to be swapped back.
Swap-Rows 9 10
Swap-Rows 29 30
Swap-Rows 49 50
Swap-Rows 76 77
Swap-Rows 87 88
Swap-Rows 99 100

# Append the new fixture as row 145, copying the row-style (bold/bordered
# index cell + date number format) from the previous last row (144).
$ws.Range("A144:AC144").Copy()
$ws.Range("A145:AC145").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A145").Value2 = 143
$ws.Range("B145").Value2 = 7952737
$ws.Range("C145").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D145").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E145").Value2 = 45381.69791666666
$ws.Range("F145").Value = "GOSK Gabela"
$ws.Range("G145").Value = "NK Posusje"
$ws.Range("H145").Value2 = 0
$ws.Range("I145").Value2 = 1
$ws.Range("J145").Value = "A"
$ws.Range("K145").Value2 = 2.5
$ws.Range("L145").Value2 = 3.2
$ws.Range("M145").Value2 = 2.5
$ws.Range("N145").Value2 = 2.7
$ws.Range("O145").Value2 = 2.8
$ws.Range("P145").Value2 = 2.7
$ws.Range("Q145").Value2 = 0
$ws.Range("R145").Value2 = 1.95
$ws.Range("S145").Value2 = 1.85
$ws.Range("T145").Value2 = 1.75
$ws.Range("U145").Value2 = 1.825
$ws.Range("V145").Value2 = 1.975
$ws.Range("W145").Value2 = -1
$ws.Range("X145").Value2 = -1
$ws.Range("Y145").Value2 = 1.7
$ws.Range("Z145").Value2 = -1
$ws.Range("AA145").Value2 = 0.8500000000000001
$ws.Range("AB145").Value2 = -1
$ws.Range("AC145").Value2 = 0.9750000000000001
